$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# Fill in the new journal entry on row 57
$ws.Range("A57").Value = Get-Date -Year 2023 -Month 6 -Day 22 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("B57").Value = "Implémentation"
$ws.Range("C57").Value = 3
$ws.Range("D57").Value = "Backend: refactor global"

# Move the active selection to D59 (next empty "Travail effectué" cell)
$ws.Range("D59").Select()
